$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TextBox 41" (shape id=42) holds two file-path lines; the script
# referenced on the second line was renamed from gdx2txt.gms to
# gdx4png.gms. Update just that run's text so its run properties
# (rPr) stay untouched.
$sh = $s.Shapes.Item("TextBox 41")
$run = $sh.TextFrame.TextRange.Paragraphs(2).Runs(1)
$run.Text = "\prog\gdx4png.gms"
